$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E4").Value = "2016-03-31 06:56:12"
$wsZh.Range("H4").Value = "2016-03-31 06:57:11"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E4").Value = "2016-03-31 06:56:22"
$wsDe.Range("H4").Value = "2016-03-31 06:57:27"
